# "pushed updates to leaflet" - refresh proposed IPTDS site data on Sheet1:
# update action_priority/notes values, append trailing periods to several
# notes, and fix site ordering/coordinates for a few proposed sites.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3
$ws.Range("F3").Value = "If transfer to IPTDS O&M project is not desired, ensure long-term funding."

# Row 4
$ws.Range("A4").Value = "USU (Proposed)"
$ws.Range("C4").Value = 44.253700000000002
$ws.Range("D4").Value = -114.35290000000001
$ws.Range("F4").Value = "Alternative locations: East Fork Salmon River."

# Row 5
$ws.Range("A5").Value = "EFS (Proposed)"
$ws.Range("C5").Value = 44.245899999999999
$ws.Range("D5").Value = -114.30289999999999
$ws.Range("F5").ClearContents()

# Row 6
$ws.Range("A6").Value = "USI (Proposed)"
$ws.Range("C6").Value = 44.677300000000002
$ws.Range("D6").Value = -114.0744
$ws.Range("F6").Value = "Ideally, located below Morgan Creek and above population boundary. Locations near Challis, ID could also be considered."

# Row 7
$ws.Range("A7").Value = "USI"
$ws.Range("B7").Value = "Decommission, Remove, or Transfer"
$ws.Range("C7").Value = 44.889763000000002
$ws.Range("D7").Value = -113.964145
$ws.Range("E7").Value = "LOW"
$ws.Range("F7").Value = "Move upstream to proposed USI location."

# Row 8
$ws.Range("A8").Value = "USE"
$ws.Range("B8").Value = "Continue Funding"
$ws.Range("C8").Value = 45.028530000000003
$ws.Range("D8").Value = -113.916319
$ws.Range("E8").Value = "MED"
$ws.Range("F8").Value = "Upgrade to IS1001 MC to span river and increase read range."

# Row 12
$ws.Range("E12").Value = "MED"
$ws.Range("F12").Value = "Planned upsgrade to IS1001 MC. Alternatively, LLR could be considered for upgrade to IS1001 MC to increase juvenile detections."

# Row 17
$ws.Range("F17").Value = "If transfer to IPTDS O&M project is not desired, ensure long-term funding."

# Row 18
$ws.Range("F18").Value = "Alternate locations: upper Chamberlain, Sabe, Bargamin, Warren, Crooked, or Sheep creeks."

# Row 23
$ws.Range("F23").Value = "Upgrade upstream and/or downstream array(s) to FS1001 MUX."

# Row 24
$ws.Range("F24").Value = "Upgrade upstream and/or downstream array(s) to FS1001 MUX."

# Row 25
$ws.Range("F25").Value = "Alternative locations: Slate Creek, Whitebird Creek."

# Row 28
$ws.Range("E28").Value = "HIGH"

# Row 29
$ws.Range("E29").Value = "HIGH"

# Row 32
$ws.Range("F32").Value = "If not already, consider upgrade to IS1001 to increase read range."

# Row 33
$ws.Range("F33").Value = "Move to lower boundary of CRSFC-s population."

# Row 34
$ws.Range("F34").Value = "Consider funding either SC3 or SC4 to provide detections to estimate abundance at SC2 at proposed new location."

# Row 35
$ws.Range("F35").Value = "Consider funding either SC3 or SC4 to provide detections to estimate abundance at SC2 at proposed new location."

# Row 37
$ws.Range("F37").Value = "Consider upgrade to IS1001 MC and tandem arrays if sufficient distance btw arrays can be achieved."

# Row 38
$ws.Range("F38").Value = "Consider decommissioning if LC1 can be converted to a two-pass configuration; alternatively, consider moving upstream to below core spawning areas."

# Row 39
$ws.Range("F39").Value = "One of LAP, JUL, or LAW should be funded to provide monitoring in the CRLMA-s population. LAP or JUL would be preferred to continue time-series of estimates."

# Row 40
$ws.Range("F40").Value = "One of LAP, JUL, or LAW should be funded to provide monitoring in the CRLMA-s population. LAP or JUL would be preferred to continue time-series of estimates."

# Row 41
$ws.Range("F41").Value = "One of LAP, JUL, or LAW should be funded to provide monitoring in the CRLMA-s population. LAP or JUL would be preferred to continue time-series of estimates."

# Row 42
$ws.Range("F42").Value = "Upgrade to IS1001 MC to increase read range."

# Row 47
$ws.Range("F47").Value = "If not already, consider upgrade to IS1001 to increase read range."

# Row 48
$ws.Range("F48").Value = "Upgrade to IS1001 MC to span river and increase read range."

# Row 49
$ws.Range("F49").Value = "New site could be a single-pass array."

# Row 50
$ws.Range("F50").Value = "If transfer to IPTDS O&M project is not desired, ensure long-term funding."

# Row 51
$ws.Range("F51").Value = "If transfer to IPTDS O&M project is not desired, ensure long-term funding."

# Row 54
$ws.Range("F54").Value = "Proposed site not necessary if weir is operated annually and is reliable."

$ws.Range("F53").Select()
Write-Output "done"
